# Corrijo datos y graficas
# Update the metrics table (rows 2-6, columns B-M) with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1695.0, 936.0, 35612.0, 2945.0, 0.6442417331812998, 0.36530172413793105, 0.974389843493488, 0.46623573098610926, 0.9057735262697874, 0.41885657686789696, 0.5966123446683449, 0.923807018781053)
    3 = @(1702.0, 939.0, 35609.0, 2938.0, 0.6444528587656191, 0.36681034482758623, 0.9743077596585312, 0.4675181980497184, 0.905870641934544, 0.42012820708539766, 0.5978178362081035, 0.9238504430167439)
    4 = @(1712.0, 932.0, 35616.0, 2928.0, 0.6475037821482602, 0.3689655172413793, 0.9744992886067637, 0.4700713893465129, 0.9062833835097601, 0.422871461331635, 0.5996304145656313, 0.923903558688107)
    5 = @(1719.0, 929.0, 35619.0, 2921.0, 0.6491691842900302, 0.3704741379310345, 0.9745813724417205, 0.4717343578485181, 0.9065261726716519, 0.4246331876831825, 0.6008803489872097, 0.9238813777723403)
    6 = @(1729.0, 926.0, 35622.0, 2911.0, 0.6512241054613936, 0.3726293103448276, 0.9746634562766773, 0.47402330363262507, 0.9068417985821113, 0.4270403235986943, 0.6026509533143413, 0.9237469772862761)
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # Column B is index 2
        $ws.Cells.Item($rowNum, $col).Value = $values[$i]
    }
}

$wb.Save()
